$wb = $excel.ActiveWorkbook

# Sheet 1: Weekly Quantity
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B3").Value = 1
$ws1.Range("A4:B7").EntireRow.Delete()

# Sheet 2: Monthly Trend
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B3").Value = 1
$ws2.Range("A4:B4").EntireRow.Delete()
